$wb = $excel.ActiveWorkbook

# Sheet1 ("SingleLogin") is the currently-active/selected sheet in the
# source workbook (tabSelected="1", selection at B20). Move the selection
# back to its default (A1) before we switch away, so it no longer carries
# a stale selection once it loses focus.
$singleLogin = $wb.Worksheets.Item(1)
$excel.Goto($singleLogin.Range("A1"), $false)

# Add the new "HomePage" worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "HomePage"

# Fill the sheet column-by-column (matches the author's shared-string order).
$ws.Range("A1").Value = "SearchText"
$ws.Range("A2").Value = "SUMMER DRESSES"
$ws.Range("B1").Value = "SearchTextResult"
$ws.Range("B2").Value = "Printed Summer Dress`n"
$ws.Range("C1").Value = "SearchText1"
$ws.Range("C2").Value = "Chiffon"

# Formatting: centre the second row's search-term/result cells, then apply
# the bold+yellow header look to row 1 (A1:B1 bold+fill reuses the existing
# header style; C1 gets fill only, no bold).
$ws.Range("A2:B2").HorizontalAlignment = -4108
$ws.Range("A1:C1").Interior.Color = 65535
$ws.Range("A1:B1").Font.Bold = $true

# HomePage becomes the active sheet/tab, with C2 selected.
$ws.Activate()
$ws.Range("C2").Select()

Write-Output "done"
